$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '71.020.06'
$ws.Range("E2").Value = '  +6.23%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.797.31'
$ws.Range("E3").Value = '  +23.25%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.23%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '613.98'
$ws.Range("E5").Value = '  +7.54%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '180.27'
$ws.Range("E6").Value = '  +2.23%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.804.20'
$ws.Range("E7").Value = '  +23.47%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.16%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.546'
$ws.Range("E9").Value = '  +6.74%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.171'
$ws.Range("E10").Value = '  +12.98%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.42'
$ws.Range("E11").Value = '  +0.48%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.505'
$ws.Range("E12").Value = '  +8.46%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '41.21'
$ws.Range("E13").Value = '  +15.31%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000262'
$ws.Range("E14").Value = '  +9.42%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.421.67'
$ws.Range("E15").Value = '  +23.03%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.789.96'
$ws.Range("E16").Value = '  +22.84%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '71.091.84'
$ws.Range("E17").Value = '  +6.28%  '

$ws.Range("E18").Value = '  +1.66%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.65'
$ws.Range("E19").Value = '  +9.32%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '526.71'
$ws.Range("E20").Value = '  +8.58%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.89'
$ws.Range("E21").Value = '  +2.68%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.51'
$ws.Range("E22").Value = '  +24.02%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.753'
$ws.Range("E23").Value = '  +10.40%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '88.93'
$ws.Range("E24").Value = '  +6.83%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.49'
$ws.Range("E25").Value = '  +11.66%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '13.71'
$ws.Range("E26").Value = '  +8.66%  '

$ws.Range("E27").Value = '  +8.20%  '

$ws.Range("E28").Value = '  +0.07%  '

$ws.Range("E29").Value = '  +32.71%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.53'
$ws.Range("E30").Value = '  +10.73%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.93'
$ws.Range("E31").Value = '  +14.00%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.08'
$ws.Range("E32").Value = '  +3.05%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '32.57'
$ws.Range("E33").Value = '  +16.74%  '

$ws.Range("E34").Value = '  +4.78%  '

$ws.Range("E35").Value = '  -0.30%  '

$ws.Range("E36").Value = '  +13.00%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.24'
$ws.Range("E37").Value = '  +12.47%  '

$ws.Range("B38").Value = 'Stacks'
$ws.Range("C38").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.24'
$ws.Range("E38").Value = '  +12.24%  '

$ws.Range("B39").Value = 'TheGraph'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.344'
$ws.Range("E39").Value = '  +11.29%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.133'
$ws.Range("E40").Value = '  +8.14%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '51.68'
$ws.Range("E41").Value = '  +5.56%  '

$ws.Range("B42").Value = 'Cosmos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.92'
$ws.Range("E42").Value = '  +8.87%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.161.83'
$ws.Range("E43").Value = '  +13.30%  '

$ws.Range("B44").Value = 'Bittensor'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '430.38'
$ws.Range("E44").Value = '  +17.20%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '44.31'
$ws.Range("E45").Value = '  -6.21%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.83'
$ws.Range("E46").Value = '  +4.50%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0371'
$ws.Range("E47").Value = '  +8.37%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '27.95'
$ws.Range("E48").Value = '  +9.92%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '141.60'
$ws.Range("E49").Value = '  +5.25%  '

$ws.Range("E50").Value = '  +0.02%  '

$ws.Range("E51").Value = '  +8.12%  '
